$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the existing "InputFile" text column
# (old D -> new E), then populate the new D column with "Step" numbers.
$ws.Columns("D:D").Insert()

# Header
$ws.Range("D1").Value = "Step"

# Fill Step values: the AFNI step numbers cycle 1,3,5,7 for each subject block.
for ($r = 2; $r -le 153; $r++) {
    $idx = ($r - 2) % 4
    if ($idx -eq 0) { $val = 1 }
    elseif ($idx -eq 1) { $val = 3 }
    elseif ($idx -eq 2) { $val = 5 }
    else { $val = 7 }
    $ws.Cells.Item($r, 4).Value = $val
}

# The inserted column should carry plain/default formatting (no inherited style).
$ws.Range("D2:D153").Style = "Normal"

# Update the current selection to match the post-edit cursor position.
$ws.Range("G35").Select() | Out-Null

# Re-apply the sort over the now-wider A:E data range so the persisted
# sort state reflects the extra column.
$ws.Sort.SetRange($ws.Range("A2:E77"))
$ws.Sort.Apply()
